# The workbook originally has two sheets, in tab order:
#   1) "2022-Q1" (quarterly holding detail)
#   2) "总计"    (summary totals)
#
# This edit re-sorts the sheet tabs so "总计" comes first, followed by
# "2022-Q1" - i.e. swap the tab order of the two existing sheets.
$wb = $excel.ActiveWorkbook

$firstSheet = $wb.Worksheets.Item(1)
$secondSheet = $wb.Worksheets.Item(2)

# Move the second sheet ("总计") so that it sits directly before the
# first sheet ("2022-Q1"), putting it into the first tab position.
$secondSheet.Move($firstSheet)
